$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.679.96"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "3.816.61"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'679.22"
$ws.Range("E5").Value = "  +8.74%  "
$ws.Range("D6").Value = "'170.22"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("D7").Value = "3.813.92"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").Value = "'7.27"
$ws.Range("E11").Value = "  +7.54%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'35.88"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "4.458.71"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "3.813.13"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "70.716.69"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "'17.69"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D21").Value = "'11.28"
$ws.Range("E21").Value = "  +18.25%  "
$ws.Range("D22").Value = "'477.73"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "'83.35"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'0.0000142"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").Value = "'12.28"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "3.966.57"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "'2.93"
$ws.Range("E31").Value = "  +10.27%  "
$ws.Range("D32").Value = "'2.30"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").Value = "'7.39"
$ws.Range("E33").Value = "  +3.72%  "
$ws.Range("D34").Value = "'29.58"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").Value = "'9.13"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "3.771.88"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "'3.38"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "'5.94"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  +11.74%  "
$ws.Range("D46").Value = "'46.14"
$ws.Range("E46").Value = "  +6.98%  "
$ws.Range("D47").Value = "'159.51"
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("D48").Value = "'48.13"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").Value = "'1.44"
$ws.Range("E49").Value = "  +6.48%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000295"
$ws.Range("E50").Value = "  +8.86%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.300"
$ws.Range("E51").Value = "  +1.97%  "
